# "Added a few more slots"
#
# This review doc had its old meta-description paragraph (bold "Meta
# description" label + plain ": Read our review..." sentence) sitting right
# under the H1 title. The edit removes that paragraph from the top and,
# at the very bottom, turns the old single italic "Prompt: ..." image-brief
# paragraph into two paragraphs: a new bold "Play Doom of Egypt for Free -
# Review 2021 | AP" line, followed by the same italic paragraph but now
# holding the former meta-description sentence (without the "Meta
# description" label/colon).

$d = $word.ActiveDocument

# --- capture the "Meta description" paragraph (2nd paragraph) before it's removed ---
$metaPara = $d.Paragraphs.Item(2)
$metaFormatted = $metaPara.Range.FormattedText

# --- split a new paragraph in front of the last paragraph, seeded with a copy ---
# --- of the meta paragraph's formatting/runs (bold label run included)        ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertAt = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertAt.FormattedText = $metaFormatted

# the new paragraph currently reads "Meta description: Read our review...";
# drop everything after the bold "Meta description" label (16 chars) ...
$newPara = $d.Paragraphs.Item($lastIndex)
$boldRunEnd = $newPara.Range.Start + 16
$tail = $d.Range($boldRunEnd, $newPara.Range.End - 1)
$tail.Delete()

# ... then retext the remaining bold run with the title line
$boldRun = $d.Range($newPara.Range.Start, $newPara.Range.Start + 16)
$boldRun.Text = "Play Doom of Egypt for Free - Review 2021 | AP"

# give the new paragraph the same leading empty run every other paragraph in
# this document uses
$newParaStart = $d.Paragraphs.Item($lastIndex).Range.Start
$leadIn = $d.Range($newParaStart, $newParaStart)
$emptyRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$leadIn.InsertXML($emptyRunXml)

# --- remove the original "Meta description" paragraph entirely ---
$d.Paragraphs.Item(2).Range.Delete()

# --- swap the old "Prompt: ..." image-brief paragraph's text for the ---
# --- meta-description sentence, keeping its existing italic run     ---
$finalIndex = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalIndex)
$finalText = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalText.Text = "Read our review of Doom of Egypt slot game and play for free. Discover what we like and don't like about the game's visuals, payouts, volatility, and theme."
